$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 34), pushing the
# existing data (previously rows 34-62) down to rows 36-64. Excel's Insert
# copies formatting (incl. the date number format on column D) from the
# row immediately above, so the new rows inherit the correct styling.
$ws.Rows("34:35").Insert()

# New record 1 (row 34) - Primera quality, week of 2022-01-06
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = "2022-01-06"
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 100112026
$ws.Range("G34").Value = "Haba"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 170
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = 15471
$ws.Range("N34").Value = "`$/saco 25 kilos"
$ws.Range("O34").Value = "Región de La Araucanía"
$ws.Range("P34").Value = 619
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"

# New record 2 (row 35) - Segunda quality, week of 2022-01-06
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = "2022-01-06"
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112026
$ws.Range("G35").Value = "Haba"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 13000
$ws.Range("N35").Value = "`$/saco 25 kilos"
$ws.Range("O35").Value = "Región de La Araucanía"
$ws.Range("P35").Value = 520
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
